$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.052.96'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.651.20'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.48'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5261'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2595'
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06317'
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.37'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07789'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.500'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '1.648.87'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5483'
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("D15").Value = '0.0₅8194'
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.43'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '26.062.98'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.575'
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.40'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.08'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.041'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.86'
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1237'
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.224'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.05'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.427'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05813'
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.271'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.547'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.261'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.585'
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("B34").Value = 'MXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.777'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.408'
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9440'
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5736'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.745'
$ws.Range("E39").Value = '  -5.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8429'
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.63'
$ws.Range("E42").Value = '  +3.23%  '
$ws.Range("D43").Value = '1.028.71'
$ws.Range("D44").Value = '1.795.40'
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '56.95'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  +3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.881'
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.468'
$ws.Range("E50").Value = '  +1.87%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.799'
$ws.Range("E51").Value = '  +1.94%  '
